$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F27").Value = -2.467098327217072
$ws.Range("F30").Value = -2.507324048495076
$ws.Range("F34").Value = -2.560511199647814
$ws.Range("F35").Value = -2.570720347186288
$ws.Range("F37").Value = -2.596795332656076
$ws.Range("F39").Value = -2.622553931646373
$ws.Range("F40").Value = -2.631175350759489
$ws.Range("F41").Value = -2.645714058450671
$ws.Range("F42").Value = -2.653887223508724
$ws.Range("F43").Value = -2.664502231941872
$ws.Range("F44").Value = -2.67764522327626
$ws.Range("F45").Value = -2.690785382253541
$ws.Range("F46").Value = -2.700221331850753
$ws.Range("F51").Value = -2.75705783333938
$ws.Range("F52").Value = -2.765703999991357
$ws.Range("F54").Value = -2.782428771427416
$ws.Range("F56").Value = -2.803262106699624
$ws.Range("F58").Value = -2.822173308484399
$ws.Range("F81").Value = -2.602455359578926
$ws.Range("F82").Value = -2.691996966314423
$ws.Range("F83").Value = -2.778118422536482
$ws.Range("F85").Value = -2.935544766773462
$ws.Range("F86").Value = -3.006312679371181
$ws.Range("F87").Value = -3.07181058190743
$ws.Range("F89").Value = -3.186705976337926
$ws.Range("F92").Value = -2.693929621677565
$ws.Range("F94").Value = -2.878502648859578
$ws.Range("F96").Value = -3.040252399997531
$ws.Range("F97").Value = -3.115174311473011
$ws.Range("F98").Value = -3.182241701641028
$ws.Range("F99").Value = -3.244330565093539
$ws.Range("F100").Value = -3.300035057240799
$ws.Range("F103").Value = -2.764657709275403
$ws.Range("F104").Value = -2.861198605218573
$ws.Range("F105").Value = -2.950219264967374
$ws.Range("F106").Value = -3.037534477934484
$ws.Range("F107").Value = -3.117176389508077
$ws.Range("F108").Value = -3.193968980215297
$ws.Range("F109").Value = -3.263275561660649
$ws.Range("F110").Value = -3.32827416359773
$ws.Range("F111").Value = -3.386076712033265
$ws.Range("F114").Value = -2.805207480416285
$ws.Range("F115").Value = -2.903127880549067
$ws.Range("F116").Value = -2.995934458656714
$ws.Range("F117").Value = -3.083819712205104
$ws.Range("F118").Value = -3.165652659545072
$ws.Range("F119").Value = -3.242360402895268
$ws.Range("F120").Value = -3.313974013823221
$ws.Range("F121").Value = -3.378322243986201
$ws.Range("F122").Value = -3.438317247717399
$ws.Range("F125").Value = -2.840504188109541
$ws.Range("F126").Value = -2.939719017880959
$ws.Range("F127").Value = -3.033335094953279
$ws.Range("F128").Value = -3.122001490590399
$ws.Range("F129").Value = -3.204449316505002
$ws.Range("F130").Value = -3.282481976532654
$ws.Range("F131").Value = -3.354681450306451
$ws.Range("F132").Value = -3.42022811029055
$ws.Range("F133").Value = -3.480943066104167
$ws.Range("F136").Value = -2.866474906327297
$ws.Range("F137").Value = -2.96469839756398
$ws.Range("F138").Value = -3.05930842439511
$ws.Range("F139").Value = -3.148298123579054
$ws.Range("F140").Value = -3.232060576305256
$ws.Range("F141").Value = -3.310862983972657
$ws.Range("F142").Value = -3.382558218351211
$ws.Range("F143").Value = -3.449284731470146
$ws.Range("F144").Value = -3.511037512684152
$ws.Range("F147").Value = -2.885815343302841
$ws.Range("F148").Value = -2.986208933037134
$ws.Range("F149").Value = -3.081492904450194
$ws.Range("F150").Value = -3.171027529204007
$ws.Range("F151").Value = -3.255589443401318
$ws.Range("F152").Value = -3.334127727679717
$ws.Range("F153").Value = -3.407265619745683
$ws.Range("F154").Value = -3.474922952384219
$ws.Range("F155").Value = -3.53659103589672
$ws.Range("F158").Value = -2.901766297685273
$ws.Range("F159").Value = -3.002460285083956
$ws.Range("F160").Value = -3.098164115908872
$ws.Range("F161").Value = -3.188141541034982
$ws.Range("F162").Value = -3.272620028315838
$ws.Range("F163").Value = -3.351869238312134
$ws.Range("F164").Value = -3.425628876886776
$ws.Range("F165").Value = -3.49324742557358
$ws.Range("F166").Value = -3.55556408635527
$ws.Range("F169").Value = -2.915690669790165
$ws.Range("F170").Value = -3.016848712998246
$ws.Range("F171").Value = -3.112522641398796
$ws.Range("F172").Value = -3.20319739182137
$ws.Range("F173").Value = -3.288085649571086
$ws.Range("F174").Value = -3.367690800262283
$ws.Range("F175").Value = -3.441173153111597
$ws.Range("F176").Value = -3.509815837756693
$ws.Range("F177").Value = -3.572610170122168
$ws.Range("F180").Value = -2.926558566883584
$ws.Range("F181").Value = -3.02808983933252
$ws.Range("F182").Value = -3.124506686776597
$ws.Range("F183").Value = -3.214786569390341
$ws.Range("F184").Value = -3.29990385762617
$ws.Range("F185").Value = -3.379832675193223
$ws.Range("F186").Value = -3.454023623978513
$ws.Range("F187").Value = -3.522875930127133
$ws.Range("F188").Value = -3.585775058106249
$ws.Range("F191").Value = -2.936656880863203
$ws.Range("F192").Value = -3.038332459194529
$ws.Range("F193").Value = -3.13657752987085
$ws.Range("F194").Value = -3.225576430629346
$ws.Range("F195").Value = -3.310875966031492
$ws.Range("F196").Value = -3.390607289008174
$ws.Range("F197").Value = -3.465131445412871
$ws.Range("F198").Value = -3.534337097027632
$ws.Range("F199").Value = -3.598008580038769
$ws.Range("F202").Value = -2.946049514524313
$ws.Range("F203").Value = -3.048574960747966
$ws.Range("F204").Value = -3.142824189401708
$ws.Range("F205").Value = -3.233757586562004
$ws.Range("F206").Value = -3.319277628599667
$ws.Range("F207").Value = -3.399466743686772
$ws.Range("F208").Value = -3.47431587721606
$ws.Range("F209").Value = -3.543779948088166
$ws.Range("F210").Value = -3.607775632287131
